$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 450, pushing existing rows 450-484 down to 453-487.
$ws.Range("A450:A452").EntireRow.Insert()

# Row 450: new "Clemenuless" / Especial entry (Región de O'Higgins)
$ws.Cells.Item(450, 1).Value2 = 9
$ws.Cells.Item(450, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(450, 3).Value2 = "Metropolitana"
$ws.Cells.Item(450, 4).Value2 = 44714
$ws.Cells.Item(450, 5).Value2 = 13
$ws.Cells.Item(450, 6).Value2 = "Fruta"
$ws.Cells.Item(450, 7).Value2 = 100102
$ws.Cells.Item(450, 8).Value2 = "Cítricos"
$ws.Cells.Item(450, 9).Value2 = 100102004
$ws.Cells.Item(450, 10).Value2 = "Mandarina"
$ws.Cells.Item(450, 11).Value2 = "Clemenuless"
$ws.Cells.Item(450, 12).Value2 = "Especial"
$ws.Cells.Item(450, 13).Value2 = 290
$ws.Cells.Item(450, 14).Value2 = 10000
$ws.Cells.Item(450, 15).Value2 = 10000
$ws.Cells.Item(450, 16).Value2 = 10000
$ws.Cells.Item(450, 17).Value2 = "`$/bandeja 10 kilos"
$ws.Cells.Item(450, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(450, 19).Value2 = 1000
$ws.Cells.Item(450, 20).Value2 = 10

# Row 451: new "Clemenuless" / Primera entry (Región de O'Higgins)
$ws.Cells.Item(451, 1).Value2 = 9
$ws.Cells.Item(451, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(451, 3).Value2 = "Metropolitana"
$ws.Cells.Item(451, 4).Value2 = 44714
$ws.Cells.Item(451, 5).Value2 = 13
$ws.Cells.Item(451, 6).Value2 = "Fruta"
$ws.Cells.Item(451, 7).Value2 = 100102
$ws.Cells.Item(451, 8).Value2 = "Cítricos"
$ws.Cells.Item(451, 9).Value2 = 100102004
$ws.Cells.Item(451, 10).Value2 = "Mandarina"
$ws.Cells.Item(451, 11).Value2 = "Clemenuless"
$ws.Cells.Item(451, 12).Value2 = "Primera"
$ws.Cells.Item(451, 13).Value2 = 300
$ws.Cells.Item(451, 14).Value2 = 8000
$ws.Cells.Item(451, 15).Value2 = 8000
$ws.Cells.Item(451, 16).Value2 = 8000
$ws.Cells.Item(451, 17).Value2 = "`$/bandeja 10 kilos"
$ws.Cells.Item(451, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(451, 19).Value2 = 800
$ws.Cells.Item(451, 20).Value2 = 10

# Row 452: new "Clemenuless" / Segunda entry (Región de O'Higgins)
$ws.Cells.Item(452, 1).Value2 = 9
$ws.Cells.Item(452, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(452, 3).Value2 = "Metropolitana"
$ws.Cells.Item(452, 4).Value2 = 44714
$ws.Cells.Item(452, 5).Value2 = 13
$ws.Cells.Item(452, 6).Value2 = "Fruta"
$ws.Cells.Item(452, 7).Value2 = 100102
$ws.Cells.Item(452, 8).Value2 = "Cítricos"
$ws.Cells.Item(452, 9).Value2 = 100102004
$ws.Cells.Item(452, 10).Value2 = "Mandarina"
$ws.Cells.Item(452, 11).Value2 = "Clemenuless"
$ws.Cells.Item(452, 12).Value2 = "Segunda"
$ws.Cells.Item(452, 13).Value2 = 280
$ws.Cells.Item(452, 14).Value2 = 6000
$ws.Cells.Item(452, 15).Value2 = 6000
$ws.Cells.Item(452, 16).Value2 = 6000
$ws.Cells.Item(452, 17).Value2 = "`$/bandeja 10 kilos"
$ws.Cells.Item(452, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(452, 19).Value2 = 600
$ws.Cells.Item(452, 20).Value2 = 10

# Row 450-452's date column (D) uses the same date-number-format style as the rest of column D.
$ws.Range("D453").Copy()
$ws.Range("D450:D452").PasteSpecial(-4122)
$excel.CutCopyMode = $false
